$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B12 was stored as inline string "3"; change it to a true numeric value 3.
$ws.Range("B12").Value = 3

# Add new row 13 with the annotation data.
$ws.Range("A13").Value = "Sunsi Wu"
$ws.Range("B13").Value = "'1"
$ws.Range("C13").Value = "absolutely"
$ws.Range("D13").Value = "CRT"
$ws.Range("E13").Value = "EXP"
$ws.Range("F13").Value = "77474e59-42ef-43e4-850b-a07d6b41a266"
$ws.Range("G13").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H13").Value = "You absolutely know this but you hide these results."
